# Atualização automática da planilha
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")
$ws.Activate()

# Advance the timeline header in Q1 by one month (2026-11-01 -> 2026-12-01)
$ws.Range("Q1").Value = 46357

# Clear the now-unneeded forward-looking monthly forecast cells (H:Q) for every budget line,
# keeping their existing formatting/styles intact.
$ws.Range("H3:Q11").ClearContents()

# Update row 4 (Consultoria & Implantação / DBD / P01): collapse forecast into current actuals
$ws.Range("F4").Value = 47200
$ws.Range("G4").Value = 47200

# Move the active selection to G4 (matches the saved selection state)
$ws.Range("G4").Select()
